# Venezuela Primera Division workbook update (19-02-2024 20:58)
# The underlying change is a re-sort of several match rows: full records
# (columns B and F:AC) moved between row positions while column A (the
# positional running index) and columns C/D/E (Div/Div Original Name/Date,
# identical across the affected rows) stayed put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot original row data (col B, and F:AC) before overwriting anything ---
$rowB_197 = $ws.Range("B197").Value()
$rowFAC_197 = $ws.Range("F197:AC197").Value()
$rowB_198 = $ws.Range("B198").Value()
$rowFAC_198 = $ws.Range("F198:AC198").Value()
$rowB_200 = $ws.Range("B200").Value()
$rowFAC_200 = $ws.Range("F200:AC200").Value()
$rowB_201 = $ws.Range("B201").Value()
$rowFAC_201 = $ws.Range("F201:AC201").Value()
$rowB_202 = $ws.Range("B202").Value()
$rowFAC_202 = $ws.Range("F202:AC202").Value()
$rowB_203 = $ws.Range("B203").Value()
$rowFAC_203 = $ws.Range("F203:AC203").Value()
$rowB_206 = $ws.Range("B206").Value()
$rowFAC_206 = $ws.Range("F206:AC206").Value()
$rowB_207 = $ws.Range("B207").Value()
$rowFAC_207 = $ws.Range("F207:AC207").Value()
$rowB_239 = $ws.Range("B239").Value()
$rowFAC_239 = $ws.Range("F239:AC239").Value()
$rowB_240 = $ws.Range("B240").Value()
$rowFAC_240 = $ws.Range("F240:AC240").Value()

# --- Step 2: write snapshots into their new (destination) row positions ---
$ws.Range("B197").Value = $rowB_198
$ws.Range("F197:AC197").Value = $rowFAC_198
$ws.Range("B198").Value = $rowB_200
$ws.Range("F198:AC198").Value = $rowFAC_200
$ws.Range("B200").Value = $rowB_203
$ws.Range("F200:AC200").Value = $rowFAC_203
$ws.Range("B201").Value = $rowB_197
$ws.Range("F201:AC201").Value = $rowFAC_197
$ws.Range("B202").Value = $rowB_201
$ws.Range("F202:AC202").Value = $rowFAC_201
$ws.Range("B203").Value = $rowB_202
$ws.Range("F203:AC203").Value = $rowFAC_202
$ws.Range("B206").Value = $rowB_207
$ws.Range("F206:AC206").Value = $rowFAC_207
$ws.Range("B207").Value = $rowB_206
$ws.Range("F207:AC207").Value = $rowFAC_206
$ws.Range("B239").Value = $rowB_240
$ws.Range("F239:AC239").Value = $rowFAC_240
$ws.Range("B240").Value = $rowB_239
$ws.Range("F240:AC240").Value = $rowFAC_239
